# Add a new "alphabet-war" exercise column (F) to the tracking sheet,
# and correct the "sql-basics-monsters-using-case" (E) results now that
# the CSV parser trims/allows spaces.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F, matching the style used by the other header cells.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "alphabet-war"

# Corrected values for column E (sql-basics-monsters-using-case) and the
# newly populated column F (alphabet-war) for every data row.
$rows = @(
    @{Row=2; E=$False; F=$True},
    @{Row=3; E=$True; F=$True},
    @{Row=4; E=$True; F=$False},
    @{Row=5; E=$True; F=$False},
    @{Row=6; E=$True; F=$False},
    @{Row=7; E=$False; F=$False},
    @{Row=8; E=$True; F=$False},
    @{Row=9; E=$True; F=$True},
    @{Row=10; E=$False; F=$False},
    @{Row=11; E=$True; F=$True},
    @{Row=12; E=$True; F=$False},
    @{Row=13; E=$True; F=$False},
    @{Row=14; E=$True; F=$False},
    @{Row=15; E=$True; F=$False},
    @{Row=16; E=$True; F=$False},
    @{Row=17; E=$False; F=$False},
    @{Row=18; E=$True; F=$False},
    @{Row=19; E=$True; F=$True},
    @{Row=20; E=$True; F=$True},
    @{Row=21; E=$True; F=$False}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
